$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.162.72"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.854.48"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.71%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "238.15"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.20%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.621"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.71%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "41.77"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +4.42%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.327"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0692"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0991"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "2.121.64"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").Value = "1.871.67"
$ws.Range("E13").Value = "  +2.79%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "11.39"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.44%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.674"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.86%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "4.72"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "35.134.96"
$ws.Range("E17").Value = "  +0.35%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "69.98"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "0.0$([char]0x2083)0792"
$ws.Range("E19").Value = "  +0.52%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "240.42"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "12.14"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("E23").Value = "  +0.63%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "169.07"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.57%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.98"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.79%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.82"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +20.26%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "17.56"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("E29").Value = "  -0.43%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.72%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0553"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.58%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.98"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.41%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.01"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +27.86%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.00"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +9.01%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.800"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +15.12%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.29"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("E38").Value = "  +8.08%  "
$ws.Range("E39").Value = "  +3.34%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "89.83"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.27%  "
$ws.Range("D41").Value = "1.342.30"
$ws.Range("E41").Value = "  +0.14%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "14.78"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.58%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "12.77"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +47.55%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.29"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0555"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +6.60%  "
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").Value = "2.029.00"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0676"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.73%  "
